$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("D1").Value = "distracted"
$ws.Range("F1").Value = "FPS"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null

# Update data rows (A path separator, C/D/E values, F new FPS column)
$ws.Range("A2").Value = 'dataset/val\videos\1080217202501_dms_drowsy_1.mp4'
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 10

$ws.Range("A3").Value = 'dataset/val\videos\1120211202501_dms_drowsy_1.mp4'
$ws.Range("C3").Value = 102
$ws.Range("D3").Value = 168
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 9

$ws.Range("A4").Value = 'dataset/val\videos\1140217202501_dms_drowsy_5.mp4'
$ws.Range("C4").Value = 12
$ws.Range("D4").Value = 8
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 7

$ws.Range("A5").Value = 'dataset/val\videos\270217202501_dms_drowsy_10.mp4'
$ws.Range("C5").Value = 82
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 10

$ws.Range("A6").Value = 'dataset/val\videos\270217202501_dms_drowsy_11.mp4'
$ws.Range("C6").Value = 144
$ws.Range("D6").Value = 142
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 9

$ws.Range("A7").Value = 'dataset/val\videos\270217202501_dms_drowsy_12.mp4'
$ws.Range("C7").Value = 144
$ws.Range("D7").Value = 142
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 9

$ws.Range("A8").Value = 'dataset/val\videos\270217202501_dms_drowsy_21.mp4'
$ws.Range("C8").Value = 83
$ws.Range("D8").Value = 202
$ws.Range("E8").Value = 17
$ws.Range("F8").Value = 9

$ws.Range("A9").Value = 'dataset/val\videos\270217202501_dms_drowsy_22.mp4'
$ws.Range("C9").Value = 169
$ws.Range("D9").Value = 119
$ws.Range("E9").Value = 19
$ws.Range("F9").Value = 10

$ws.Range("A10").Value = 'dataset/val\videos\270217202501_dms_drowsy_24.mp4'
$ws.Range("C10").Value = 185
$ws.Range("D10").Value = 98
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 10

$ws.Range("A11").Value = 'dataset/val\videos\270217202501_dms_drowsy_27.mp4'
$ws.Range("C11").Value = 79
$ws.Range("D11").Value = 163
$ws.Range("E11").Value = 19
$ws.Range("F11").Value = 10

$ws.Range("A12").Value = 'dataset/val\videos\270217202501_dms_drowsy_28.mp4'
$ws.Range("C12").Value = 31
$ws.Range("D12").Value = 248
$ws.Range("E12").Value = 21
$ws.Range("F12").Value = 9

$ws.Range("A13").Value = 'dataset/val\videos\270217202501_dms_drowsy_29.mp4'
$ws.Range("C13").Value = 173
$ws.Range("D13").Value = 77
$ws.Range("E13").Value = 22
$ws.Range("F13").Value = 9

$ws.Range("A14").Value = 'dataset/val\videos\270217202501_dms_drowsy_6.mp4'
$ws.Range("C14").Value = 270
$ws.Range("D14").Value = 32
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 10

$ws.Range("A15").Value = 'dataset/val\videos\270217202501_dms_drowsy_7.mp4'
$ws.Range("C15").Value = 281
$ws.Range("D15").Value = 26
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 10

$ws.Range("A16").Value = 'dataset/val\videos\270217202501_dms_drowsy_8.mp4'
$ws.Range("C16").Value = 294
$ws.Range("D16").Value = 16
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 10

$ws.Range("A17").Value = 'dataset/val\videos\270217202501_dms_drowsy_9.mp4'
$ws.Range("C17").Value = 294
$ws.Range("D17").Value = 16
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 10

$ws.Range("A18").Value = 'dataset/val\videos\3926989_dms_drowsy_1.mp4'
$ws.Range("C18").Value = 268
$ws.Range("D18").Value = 13
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 9

$ws.Range("A19").Value = 'dataset/val\videos\3927679_dms_drowsy_1.mp4'
$ws.Range("C19").Value = 242
$ws.Range("D19").Value = 44
$ws.Range("E19").Value = 7
$ws.Range("F19").Value = 9

$ws.Range("A20").Value = 'dataset/val\videos\410217202501_dms_drowsy_11.mp4'
$ws.Range("C20").Value = 193
$ws.Range("D20").Value = 123
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 10

$ws.Range("A21").Value = 'dataset/val\videos\410217202501_dms_drowsy_12.mp4'
$ws.Range("C21").Value = 212
$ws.Range("D21").Value = 77
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 10

$ws.Range("A22").Value = 'dataset/val\videos\410217202501_dms_drowsy_2.mp4'
$ws.Range("C22").Value = 311
$ws.Range("D22").Value = 15
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 10

$ws.Range("A23").Value = 'dataset/val\videos\410217202501_dms_drowsy_3.mp4'
$ws.Range("C23").Value = 262
$ws.Range("D23").Value = 21
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 11

$ws.Range("A24").Value = 'dataset/val\videos\410217202501_dms_drowsy_4.mp4'
$ws.Range("C24").Value = 240
$ws.Range("D24").Value = 78
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 11

$ws.Range("A25").Value = 'dataset/val\videos\410217202501_dms_drowsy_5.mp4'
$ws.Range("C25").Value = 106
$ws.Range("D25").Value = 194
$ws.Range("E25").Value = 4
$ws.Range("F25").Value = 9

$ws.Range("A26").Value = 'dataset/val\videos\410217202501_dms_drowsy_6.mp4'
$ws.Range("C26").Value = 291
$ws.Range("D26").Value = 26
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 10

$ws.Range("A27").Value = 'dataset/val\videos\410217202501_dms_drowsy_7.mp4'
$ws.Range("C27").Value = 126
$ws.Range("D27").Value = 161
$ws.Range("E27").Value = 8
$ws.Range("F27").Value = 9

$ws.Range("A28").Value = 'dataset/val\videos\410217202501_dms_drowsy_9.mp4'
$ws.Range("C28").Value = 109
$ws.Range("D28").Value = 220
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 10

$ws.Range("A29").Value = 'dataset/val\videos\410218202502_dms_drowsy_21.mp4'
$ws.Range("C29").Value = 323
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 10

$ws.Range("A30").Value = 'dataset/val\videos\410224202501_dms_drowsy_11.mp4'
$ws.Range("C30").Value = 251
$ws.Range("D30").Value = 38
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 10

$ws.Range("A31").Value = 'dataset/val\videos\410224202501_dms_drowsy_12.mp4'
$ws.Range("C31").Value = 251
$ws.Range("D31").Value = 38
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 10

$ws.Range("A32").Value = 'dataset/val\videos\410224202501_dms_drowsy_13.mp4'
$ws.Range("C32").Value = 201
$ws.Range("D32").Value = 119
$ws.Range("E32").Value = 1
$ws.Range("F32").Value = 9

$ws.Range("A33").Value = 'dataset/val\videos\410224202501_dms_drowsy_14.mp4'
$ws.Range("C33").Value = 213
$ws.Range("D33").Value = 128
$ws.Range("E33").Value = 2
$ws.Range("F33").Value = 10

$ws.Range("A34").Value = 'dataset/val\videos\410224202501_dms_drowsy_15.mp4'
$ws.Range("C34").Value = 143
$ws.Range("D34").Value = 97
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 11

$ws.Range("A35").Value = 'dataset/val\videos\410224202501_dms_drowsy_16.mp4'
$ws.Range("C35").Value = 143
$ws.Range("D35").Value = 97
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = 11

$ws.Range("A36").Value = 'dataset/val\videos\410224202501_dms_drowsy_2.mp4'
$ws.Range("C36").Value = 159
$ws.Range("D36").Value = 172
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 10

$ws.Range("A37").Value = 'dataset/val\videos\410224202501_dms_drowsy_3.mp4'
$ws.Range("C37").Value = 310
$ws.Range("D37").Value = 14
$ws.Range("E37").Value = 2
$ws.Range("F37").Value = 11

$ws.Range("A38").Value = 'dataset/val\videos\410224202501_dms_drowsy_4.mp4'
$ws.Range("C38").Value = 269
$ws.Range("D38").Value = 13
$ws.Range("E38").Value = 1
$ws.Range("F38").Value = 10
